# Update the "想去人数" (want-to-go count) column F values on the two
# sheets that contain the full exhibition listing data: "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F-column value for most rows (shared by both sheets)
$updates = @{
    2  = 340
    4  = 10530
    6  = 957
    7  = 68
    11 = 453
    13 = 210
    15 = 3229
    18 = 709
    20 = 1047
    22 = 89
    23 = 1675
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# Row 9's F value differs slightly between the two sheets in the original
# data (7626 vs 7627), so handle it per-sheet explicitly.
$wb.Worksheets.Item("展览").Range("F9").Value = 7733
$wb.Worksheets.Item("全部类型").Range("F9").Value = 7734
